# Add team record (Wins / Losses / Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new headers, formatted like the existing header cells (e.g. A1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-43: every row gets the same team record - Wins=90, Losses=72, Ties=0.
$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 90
    $ws.Cells.Item($r, 31).Value = 72
    $ws.Cells.Item($r, 32).Value = 0
}

$wb.Save()
